$wb = $excel.ActiveWorkbook

# --- Sheet 1 ---
$ws = $wb.Worksheets.Item(1)
$data1 = @(
    @(0,43),
    @(1,45),
    @(2,46),
    @(3,47),
    @(4,49),
    @(5,50),
    @(6,52),
    @(7,53),
    @(8,55),
    @(9,56),
    @(10,58),
    @(11,59),
    @(12,61),
    @(13,62),
    @(14,64),
    @(15,66),
    @(16,67),
    @(17,69),
    @(18,71),
    @(19,73),
    @(20,75),
    @(21,77),
    @(22,79),
    @(23,81),
    @(24,83),
    @(25,85),
    @(26,87),
    @(27,90),
    @(28,92),
    @(29,95),
    @(30,98),
    @(31,101),
    @(32,104),
    @(33,108),
    @(34,112),
    @(35,116),
    @(36,122),
    @(37,130),
    @(38,130)
)
for ($i = 0; $i -lt $data1.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data1[$i][0]
    $ws.Cells.Item($r, 2).Value = $data1[$i][1]
}

# --- Sheet 2 ---
$ws = $wb.Worksheets.Item(2)
$data2 = @(
    @(0,43),
    @(1,44),
    @(2,46),
    @(3,47),
    @(4,48),
    @(5,50),
    @(6,51),
    @(7,52),
    @(8,54),
    @(9,55),
    @(10,57),
    @(11,58),
    @(12,60),
    @(13,61),
    @(14,63),
    @(15,65),
    @(16,66),
    @(17,68),
    @(18,70),
    @(19,72),
    @(20,73),
    @(21,75),
    @(22,77),
    @(23,79),
    @(24,81),
    @(25,83),
    @(26,86),
    @(27,88),
    @(28,91),
    @(29,93),
    @(30,96),
    @(31,99),
    @(32,103),
    @(33,106),
    @(34,111),
    @(35,116),
    @(36,126),
    @(37,130),
    @(38,130)
)
for ($i = 0; $i -lt $data2.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data2[$i][0]
    $ws.Cells.Item($r, 2).Value = $data2[$i][1]
}

# --- Sheet 3 ---
$ws = $wb.Worksheets.Item(3)
$data3 = @(
    @(0,43),
    @(1,44),
    @(2,46),
    @(3,47),
    @(4,48),
    @(5,49),
    @(6,51),
    @(7,52),
    @(8,53),
    @(9,55),
    @(10,56),
    @(11,58),
    @(12,59),
    @(13,61),
    @(14,62),
    @(15,64),
    @(16,65),
    @(17,67),
    @(18,68),
    @(19,70),
    @(20,72),
    @(21,74),
    @(22,75),
    @(23,77),
    @(24,79),
    @(25,81),
    @(26,83),
    @(27,86),
    @(28,88),
    @(29,90),
    @(30,93),
    @(31,96),
    @(32,99),
    @(33,102),
    @(34,106),
    @(35,111),
    @(36,118),
    @(37,127),
    @(38,127)
)
for ($i = 0; $i -lt $data3.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data3[$i][0]
    $ws.Cells.Item($r, 2).Value = $data3[$i][1]
}

# --- Sheet 4 ---
$ws = $wb.Worksheets.Item(4)
$data4 = @(
    @(0,44),
    @(1,45),
    @(2,47),
    @(3,48),
    @(4,49),
    @(5,50),
    @(6,51),
    @(7,53),
    @(8,54),
    @(9,55),
    @(10,57),
    @(11,58),
    @(12,59),
    @(13,61),
    @(14,62),
    @(15,64),
    @(16,65),
    @(17,67),
    @(18,68),
    @(19,70),
    @(20,71),
    @(21,73),
    @(22,75),
    @(23,77),
    @(24,78),
    @(25,80),
    @(26,82),
    @(27,84),
    @(28,86),
    @(29,89),
    @(30,91),
    @(31,94),
    @(32,96),
    @(33,99),
    @(34,103),
    @(35,107),
    @(36,112),
    @(37,121),
    @(38,125)
)
for ($i = 0; $i -lt $data4.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data4[$i][0]
    $ws.Cells.Item($r, 2).Value = $data4[$i][1]
}

# --- Sheet 5 ---
$ws = $wb.Worksheets.Item(5)
$data5 = @(
    @(0,46),
    @(1,47),
    @(2,48),
    @(3,50),
    @(4,51),
    @(5,52),
    @(6,53),
    @(7,54),
    @(8,56),
    @(9,57),
    @(10,58),
    @(11,59),
    @(12,61),
    @(13,62),
    @(14,63),
    @(15,65),
    @(16,66),
    @(17,67),
    @(18,69),
    @(19,70),
    @(20,72),
    @(21,73),
    @(22,75),
    @(23,77),
    @(24,78),
    @(25,80),
    @(26,82),
    @(27,84),
    @(28,86),
    @(29,88),
    @(30,90),
    @(31,93),
    @(32,95),
    @(33,98),
    @(34,101),
    @(35,104),
    @(36,108),
    @(37,114),
    @(38,124)
)
for ($i = 0; $i -lt $data5.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data5[$i][0]
    $ws.Cells.Item($r, 2).Value = $data5[$i][1]
}

# --- Sheet 6 ---
$ws = $wb.Worksheets.Item(6)
$data6 = @(
    @(0,51),
    @(1,52),
    @(2,53),
    @(3,54),
    @(4,55),
    @(5,56),
    @(6,57),
    @(7,59),
    @(8,60),
    @(9,61),
    @(10,62),
    @(11,63),
    @(12,64),
    @(13,66),
    @(14,67),
    @(15,68),
    @(16,70),
    @(17,71),
    @(18,72),
    @(19,74),
    @(20,75),
    @(21,76),
    @(22,78),
    @(23,79),
    @(24,81),
    @(25,83),
    @(26,84),
    @(27,86),
    @(28,88),
    @(29,90),
    @(30,92),
    @(31,94),
    @(32,96),
    @(33,99),
    @(34,101),
    @(35,104),
    @(36,108),
    @(37,112),
    @(38,119)
)
for ($i = 0; $i -lt $data6.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data6[$i][0]
    $ws.Cells.Item($r, 2).Value = $data6[$i][1]
}
